# 5.1.xlsx — "New Staging Branch" edit
#
# Summary of the change (per the OOXML diff):
#   * "Sheet1" is renamed to "5.1 MSQ" and becomes the active sheet; it
#     gets a single value, A1 = 5.
#   * On the "5.1" sheet, a new merged header cell F3:H3 is added holding
#     the shared string "Main Story Quest" (centered, same style as the
#     B2:E24 body block).
#   * The big title text in the merged F1:J2 cell is re-typed with the
#     spacing tidied up ("Vows of Virtue, Deeds of Cruelty (5.1)" ->
#     "Vows of Virtue,Deeds of Cruelty(5.1)") and the title is bolded.

$wb = $excel.ActiveWorkbook

$quest = $wb.Worksheets.Item("5.1")
$msq   = $wb.Worksheets.Item("Sheet1")

# --- "Sheet1" -> "5.1 MSQ" ------------------------------------------------
$msq.Name = "5.1 MSQ"
$msq.Range("A1").Value = 5

# --- New "Main Story Quest" header on the "5.1" sheet ---------------------
$quest.Range("F3:H3").HorizontalAlignment = -4108   # xlCenter
$quest.Range("F3").Value = "Main Story Quest"
$quest.Range("F3:H3").Merge()

# --- Tidy up + bold the existing title in F1:J2 ----------------------------
$quest.Range("F1").Value = "Vows of Virtue,Deeds of Cruelty(5.1)"
$quest.Range("F1:J2").Font.Bold = $true

# --- Make "5.1 MSQ" the active/selected tab --------------------------------
$msq.Activate()
